# Apply scheduled-runner market/price data refresh to the Chocobo Profits workbook.
# Each Leve table sheet (one per crafting class) has its market-price derived
# columns (H:N) refreshed with newly retrieved values.
$wb = $excel.ActiveWorkbook

# ----- ALC sheet -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 266.58823  # was 251.89473
$ws.Range("I33").Value = 286.16666  # was 246.6875
$ws.Range("J33").Value = 219.6  # was 279.66666
$ws.Range("K33").Value = 286.16666  # was 246.6875
$ws.Range("L33").Value = 219.6  # was 279.66666
$ws.Range("M33").Value = -57.16665999999998  # was -17.6875
$ws.Range("N33").Value = -677.6  # was -737.66666
$ws.Range("H115").Value = 1448.5714  # was 1501.1111
$ws.Range("I115").Value = 1448.5714  # was 1501.1111
$ws.Range("K115").Value = 4345.7142  # was 4503.3333
$ws.Range("M115").Value = -2778.7142  # was -2936.3333
$ws.Range("H137").Value = 1538624.1  # was 1589886.6
$ws.Range("I137").Value = 2802298  # was 2977394.8
$ws.Range("K137").Value = 8406894  # was 8932184.399999999
$ws.Range("M137").Value = -8404344  # was -8929634.399999999

# ----- ARM sheet -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1430.9  # was 1360.9
$ws.Range("I2").Value = 1123.8  # was 1086.5
$ws.Range("J2").Value = 1738  # was 1772.5
$ws.Range("K2").Value = 1123.8  # was 1086.5
$ws.Range("L2").Value = 1738  # was 1772.5
$ws.Range("M2").Value = -1010.8  # was -973.5
$ws.Range("N2").Value = -1964  # was -1998.5
$ws.Range("H74").Value = 5981.95  # was 7310.875
$ws.Range("I74").Value = 7143.154  # was 9039.4
$ws.Range("J74").Value = 3825.4285  # was 4430
$ws.Range("K74").Value = 7143.154  # was 9039.4
$ws.Range("L74").Value = 3825.4285  # was 4430
$ws.Range("M74").Value = -6269.154  # was -8165.4
$ws.Range("N74").Value = -5573.4285  # was -6178
$ws.Range("H77").Value = 5981.95  # was 7310.875
$ws.Range("I77").Value = 7143.154  # was 9039.4
$ws.Range("J77").Value = 3825.4285  # was 4430
$ws.Range("K77").Value = 35715.77  # was 45197
$ws.Range("L77").Value = 19127.1425  # was 22150
$ws.Range("M77").Value = -31347.77  # was -40829
$ws.Range("N77").Value = -27863.1425  # was -30886
$ws.Range("H102").Value = 3602.5  # was 4650.7144
$ws.Range("I102").Value = 3602.5  # was 4600
$ws.Range("J102").Value = 0  # was 4777.5
$ws.Range("K102").Value = 3602.5  # was 4600
$ws.Range("L102").Value = 0  # was 4777.5
$ws.Range("M102").ClearContents()  # was -2978
$ws.Range("N102").Value = -1980.5  # was -8021.5
$ws.Range("H116").Value = 1430.9  # was 1360.9
$ws.Range("I116").Value = 1123.8  # was 1086.5
$ws.Range("J116").Value = 1738  # was 1772.5
$ws.Range("K116").Value = 1123.8  # was 1086.5
$ws.Range("L116").Value = 1738  # was 1772.5
$ws.Range("M116").Value = 1170.2  # was 1207.5
$ws.Range("N116").Value = -6326  # was -6360.5

# ----- BSM sheet -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1430.9  # was 1360.9
$ws.Range("I3").Value = 1123.8  # was 1086.5
$ws.Range("J3").Value = 1738  # was 1772.5
$ws.Range("K3").Value = 1123.8  # was 1086.5
$ws.Range("L3").Value = 1738  # was 1772.5
$ws.Range("M3").Value = -1009.8  # was -972.5
$ws.Range("N3").Value = -1966  # was -2000.5
$ws.Range("H45").Value = 30000  # was 25000
$ws.Range("J45").Value = 30000  # was 25000
$ws.Range("L45").Value = 30000  # was 25000
$ws.Range("N45").Value = -31616  # was -26616
$ws.Range("H105").Value = 4903581  # was 5052131.5
$ws.Range("I105").Value = 5129827.5  # was 5292634.5
$ws.Range("K105").Value = 5129827.5  # was 5292634.5
$ws.Range("M105").Value = -5128080.5  # was -5290887.5

# ----- CRP sheet -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 720.8  # was 714.0952
$ws.Range("I107").Value = 629.0909  # was 625
$ws.Range("K107").Value = 629.0909  # was 625
$ws.Range("M107").Value = 1290.9091  # was 1295
$ws.Range("H132").Value = 3879.8  # was 4005.375
$ws.Range("I132").Value = 4096.3335  # was 4100.3335
$ws.Range("J132").Value = 3679.923  # was 3910.4167
$ws.Range("K132").Value = 12289.0005  # was 12301.0005
$ws.Range("L132").Value = 11039.769  # was 11731.2501
$ws.Range("M132").Value = -9759.000499999998  # was -9771.000499999998
$ws.Range("N132").Value = -16099.769  # was -16791.2501
$ws.Range("H134").Value = 2016.2106  # was 2558.4285
$ws.Range("I134").Value = 1020.5333  # was 1281.8
$ws.Range("K134").Value = 3061.5999  # was 3845.4
$ws.Range("M134").Value = -526.5999000000002  # was -1310.4

# ----- CUL sheet -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 83773.414  # was 91361
$ws.Range("I107").Value = 430  # was 400
$ws.Range("J107").Value = 143304.42  # was 200514.2
$ws.Range("K107").Value = 1290  # was 1200
$ws.Range("L107").Value = 429913.26  # was 601542.6000000001
$ws.Range("M107").Value = 630  # was 720
$ws.Range("N107").Value = -433753.26  # was -605382.6000000001
$ws.Range("H113").Value = 6250660  # was 5208996
$ws.Range("I113").Value = 699.5  # was 696
$ws.Range("J113").Value = 12500620  # was 10417296
$ws.Range("K113").Value = 2098.5  # was 2088
$ws.Range("L113").Value = 37501860  # was 31251888
$ws.Range("M113").Value = 71.5  # was 82
$ws.Range("N113").Value = -37506200  # was -31256228
$ws.Range("H122").Value = 2585.7708  # was 2350.3962
$ws.Range("I122").Value = 837.75  # was 864.9286
$ws.Range("J122").Value = 3459.7812  # was 2883.641
$ws.Range("K122").Value = 7539.75  # was 7784.3574
$ws.Range("L122").Value = 31138.0308  # was 25952.769
$ws.Range("M122").Value = -5089.75  # was -5334.3574
$ws.Range("N122").Value = -36038.0308  # was -30852.769
$ws.Range("H132").Value = 2282.6667  # was 2364.5881
$ws.Range("I132").Value = 991.53845  # was 1000
$ws.Range("K132").Value = 8923.84605  # was 9000
$ws.Range("M132").Value = -6393.84605  # was -6470

# ----- GSM sheet -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 125001550  # was 50002624
$ws.Range("I80").Value = 125001550  # was 83335370
$ws.Range("J80").Value = 0  # was 3503
$ws.Range("K80").Value = 125001550  # was 83335370
$ws.Range("L80").Value = 0  # was 3503
$ws.Range("M80").ClearContents()  # was -83334372
$ws.Range("N80").Value = -125000552  # was -5499
$ws.Range("H83").Value = 125001550  # was 50002624
$ws.Range("I83").Value = 125001550  # was 83335370
$ws.Range("J83").Value = 0  # was 3503
$ws.Range("K83").Value = 625007750  # was 416676850
$ws.Range("L83").Value = 0  # was 17515
$ws.Range("M83").ClearContents()  # was -416671858
$ws.Range("N83").Value = -625002758  # was -27499
$ws.Range("H126").Value = 3037.77  # was 3246.14
$ws.Range("I126").Value = 2865.9126  # was 2961.3333
$ws.Range("J126").Value = 3725.2  # was 4100.56
$ws.Range("K126").Value = 8597.737800000001  # was 8883.999899999999
$ws.Range("L126").Value = 11175.6  # was 12301.68
$ws.Range("M126").Value = -6127.737800000001  # was -6413.999899999999
$ws.Range("N126").Value = -16115.6  # was -17241.68
$ws.Range("H133").Value = 41653.332  # was 41746.668
$ws.Range("J133").Value = 41653.332  # was 41746.668
$ws.Range("L133").Value = 41653.332  # was 41746.668
$ws.Range("N133").Value = -51773.332  # was -51866.668

# ----- LTW sheet -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 216.31818  # was 230.38889
$ws.Range("I55").Value = 173.3077  # was 180.3
$ws.Range("J55").Value = 278.44446  # was 293
$ws.Range("K55").Value = 173.3077  # was 180.3
$ws.Range("L55").Value = 278.44446  # was 293
$ws.Range("M55").Value = -0.3077000000000112  # was -7.300000000000011
$ws.Range("N55").Value = -624.4444599999999  # was -639
$ws.Range("H68").Value = 647.3  # was 647.95
$ws.Range("I68").Value = 647.3  # was 647.95
$ws.Range("K68").Value = 647.3  # was 647.95
$ws.Range("M68").Value = 101.7  # was 101.05
$ws.Range("H71").Value = 647.3  # was 647.95
$ws.Range("I71").Value = 647.3  # was 647.95
$ws.Range("K71").Value = 3236.5  # was 3239.75
$ws.Range("M71").Value = 507.5  # was 504.25
$ws.Range("H82").Value = 1877.5927  # was 1831.3214
$ws.Range("I82").Value = 696.0833  # was 687.3077
$ws.Range("K82").Value = 696.0833  # was 687.3077
$ws.Range("M82").Value = -335.0833  # was -326.3077
$ws.Range("H85").Value = 1877.5927  # was 1831.3214
$ws.Range("I85").Value = 696.0833  # was 687.3077
$ws.Range("K85").Value = 696.0833  # was 687.3077
$ws.Range("M85").Value = 551.9167  # was 560.6923
$ws.Range("H127").Value = 29916.072  # was 31470
$ws.Range("J127").Value = 29916.072  # was 31470
$ws.Range("L127").Value = 29916.072  # was 31470
$ws.Range("N127").Value = -39836.072  # was -41390

# ----- WVR sheet -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 154000.38  # was 38333.332
$ws.Range("I62").Value = 25000  # was 8000
$ws.Range("J62").Value = 369001  # was 53500
$ws.Range("K62").Value = 25000  # was 8000
$ws.Range("L62").Value = 369001  # was 53500
$ws.Range("M62").Value = -24376  # was -7376
$ws.Range("N62").Value = -370249  # was -54748
$ws.Range("H65").Value = 154000.38  # was 38333.332
$ws.Range("I65").Value = 25000  # was 8000
$ws.Range("J65").Value = 369001  # was 53500
$ws.Range("K65").Value = 125000  # was 40000
$ws.Range("L65").Value = 1845005  # was 267500
$ws.Range("M65").Value = -121880  # was -36880
$ws.Range("N65").Value = -1851245  # was -273740
$ws.Range("H81").Value = 6000  # was 2560
$ws.Range("I81").Value = 2000  # was 2750
$ws.Range("J81").Value = 10000  # was 1800
$ws.Range("K81").Value = 4000  # was 5500
$ws.Range("L81").Value = 20000  # was 3600
$ws.Range("M81").Value = -2939  # was -4439
$ws.Range("N81").Value = -22122  # was -5722
$ws.Range("H84").Value = 6000  # was 2560
$ws.Range("I84").Value = 2000  # was 2750
$ws.Range("J84").Value = 10000  # was 1800
$ws.Range("K84").Value = 20000  # was 27500
$ws.Range("L84").Value = 100000  # was 18000
$ws.Range("M84").Value = -14696  # was -22196
$ws.Range("N84").Value = -110608  # was -28608
$ws.Range("H107").Value = 215.5  # was 224.6
$ws.Range("I107").Value = 212.66667  # was 209.42857
$ws.Range("J107").Value = 219.75  # was 260
$ws.Range("K107").Value = 638.00001  # was 628.28571
$ws.Range("L107").Value = 659.25  # was 780
$ws.Range("M107").Value = 1281.99999  # was 1291.71429
$ws.Range("N107").Value = -4499.25  # was -4620
$ws.Range("H109").Value = 27777  # was 27776.5
$ws.Range("J109").Value = 27777  # was 27776.5
$ws.Range("L109").Value = 27777  # was 27776.5
$ws.Range("N109").Value = -30551  # was -30550.5
$ws.Range("H113").Value = 373.57144  # was 352
$ws.Range("J113").Value = 373.57144  # was 352
$ws.Range("L113").Value = 1120.71432  # was 1056
$ws.Range("N113").Value = -5460.71432  # was -5396
